$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tareas")

# Assign responsables (new shared-string values) to previously-empty cells
$ws.Range("B7").Value = "Abel"
$ws.Range("B17").Value = "Silverio"
$ws.Range("B19").Value = "Silverio"
$ws.Range("B21").Value = "Silverio"

# Update the view: scroll back to top and move selection to A13
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A13").Select()
